$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 9000
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H70").Value = 2058.7273
$ws.Range("I70").Value = 1963
$ws.Range("K70").Value = 5889
$ws.Range("M70").Value = -5619
$ws.Range("H72").Value = 9000
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H73").Value = 2058.7273
$ws.Range("I73").Value = 1963
$ws.Range("K73").Value = 5889
$ws.Range("M73").Value = -4953
$ws.Range("H92").Value = 3402133.8
$ws.Range("I92").Value = 562
$ws.Range("K92").Value = 562
$ws.Range("M92").Value = 686
$ws.Range("H101").Value = 950.375
$ws.Range("I101").Value = 523
$ws.Range("J101").Value = 1499.8572
$ws.Range("K101").Value = 1569
$ws.Range("L101").Value = 4499.571599999999
$ws.Range("M101").Value = 53
$ws.Range("N101").Value = -7743.571599999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 71620
$ws.Range("J123").Value = 71620
$ws.Range("L123").Value = 71620
$ws.Range("N123").Value = -81420
$ws.Range("H132").Value = 11890.719
$ws.Range("I132").Value = 5789.107
$ws.Range("K132").Value = 17367.321
$ws.Range("M132").Value = -14837.321
$ws.Range("H133").Value = 99994.75
$ws.Range("J133").Value = 99994.75
$ws.Range("L133").Value = 99994.75
$ws.Range("N133").Value = -105054.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3058.4546
$ws.Range("I86").Value = 2527
$ws.Range("K86").Value = 2527
$ws.Range("M86").Value = -1404
$ws.Range("H89").Value = 3058.4546
$ws.Range("I89").Value = 2527
$ws.Range("K89").Value = 12635
$ws.Range("M89").Value = -7019
$ws.Range("H94").Value = 1133
$ws.Range("I94").Value = 803.9583
$ws.Range("K94").Value = 803.9583
$ws.Range("M94").Value = -352.9583
$ws.Range("H107").Value = 1357.65
$ws.Range("I107").Value = 1537.7142
$ws.Range("K107").Value = 1537.7142
$ws.Range("M107").Value = 382.2858000000001
$ws.Range("H134").Value = 3077.2646
$ws.Range("I134").Value = 3077.2646
$ws.Range("K134").Value = 9231.793799999999
$ws.Range("M134").Value = -6696.793799999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1041.9474
$ws.Range("I22").Value = 240.625
$ws.Range("J22").Value = 1624.7273
$ws.Range("K22").Value = 240.625
$ws.Range("L22").Value = 1624.7273
$ws.Range("M22").Value = 109.375
$ws.Range("N22").Value = -2324.7273
$ws.Range("H97").Value = 54185
$ws.Range("J97").Value = 53500
$ws.Range("L97").Value = 53500
$ws.Range("N97").Value = -55482
$ws.Range("H99").Value = 4027.9
$ws.Range("I99").Value = 2763.1667
$ws.Range("J99").Value = 5925
$ws.Range("K99").Value = 2763.1667
$ws.Range("L99").Value = 5925
$ws.Range("M99").Value = -1265.1667
$ws.Range("N99").Value = -8921
$ws.Range("H126").Value = 4027.9
$ws.Range("I126").Value = 2763.1667
$ws.Range("J126").Value = 5925
$ws.Range("K126").Value = 8289.500100000001
$ws.Range("L126").Value = 17775
$ws.Range("M126").Value = -5819.500100000001
$ws.Range("N126").Value = -22715
$ws.Range("H134").Value = 896.4
$ws.Range("I134").Value = 812.9167
$ws.Range("K134").Value = 2438.7501
$ws.Range("M134").Value = 96.2498999999998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 276.22223
$ws.Range("I97").Value = 214.14285
$ws.Range("K97").Value = 642.4285500000001
$ws.Range("M97").Value = -146.4285500000001
$ws.Range("H98").Value = 507.66666
$ws.Range("I98").Value = 386.5
$ws.Range("J98").Value = 750
$ws.Range("K98").Value = 1159.5
$ws.Range("L98").Value = 2250
$ws.Range("M98").Value = 338.5
$ws.Range("N98").Value = -5246

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 99326
$ws.Range("J51").Value = 99326
$ws.Range("L51").Value = 99326
$ws.Range("N51").Value = -100344
$ws.Range("H102").Value = 1773.4
$ws.Range("I102").Value = 1379.4117
$ws.Range("J102").Value = 2610.625
$ws.Range("K102").Value = 1379.4117
$ws.Range("L102").Value = 2610.625
$ws.Range("M102").Value = 242.5882999999999
$ws.Range("N102").Value = -5854.625
$ws.Range("H107").Value = 1311.4375
$ws.Range("I107").Value = 1829.8
$ws.Range("J107").Value = 447.5
$ws.Range("K107").Value = 1829.8
$ws.Range("L107").Value = 447.5
$ws.Range("M107").Value = 90.20000000000005
$ws.Range("N107").Value = -4287.5
$ws.Range("H122").Value = 2498.5
$ws.Range("I122").Value = 2498.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7495.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5045.5
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 4550.5625
$ws.Range("I126").Value = 4550.5625
$ws.Range("K126").Value = 13651.6875
$ws.Range("M126").Value = -11181.6875

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6283.8667
$ws.Range("I136").Value = 5730
$ws.Range("K136").Value = 17190
$ws.Range("M136").Value = -14640

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 20749.75
$ws.Range("I38").Value = 20999.834
$ws.Range("J38").Value = 19999.5
$ws.Range("K38").Value = 20999.834
$ws.Range("L38").Value = 19999.5
$ws.Range("M38").Value = -20526.834
$ws.Range("N38").Value = -20945.5
$ws.Range("H96").Value = 18520018
$ws.Range("I96").Value = 18520018
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 18520018
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -18518645
$ws.Range("N96").ClearContents()
$ws.Range("H107").Value = 374.08
$ws.Range("J107").Value = 398.125
$ws.Range("L107").Value = 1194.375
$ws.Range("N107").Value = -5034.375
$ws.Range("H119").Value = 40000
$ws.Range("I119").Value = 40000
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 40000
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = -35162
$ws.Range("N119").ClearContents()
$ws.Range("H122").Value = 4203.8647
$ws.Range("I122").Value = 4421.0415
$ws.Range("J122").Value = 3802.923
$ws.Range("K122").Value = 13263.1245
$ws.Range("L122").Value = 11408.769
$ws.Range("M122").Value = -10813.1245
$ws.Range("N122").Value = -16308.769
$ws.Range("H132").Value = 1693.6571
$ws.Range("I132").Value = 1693.6571
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5080.971299999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2550.971299999999
$ws.Range("N132").ClearContents()
